$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A "Deposit" of 2 BTC to Poloniex (row 52) was missing its matching
# "Withdrawal" from the source wallet, which failed the integrity check.
# Insert a new row above the existing Deposit row to record that
# Withdrawal, shifting all the following transaction rows down by one.
$ws.Rows("52:52").Insert()

$ws.Range("A52").Value2 = "Withdrawal"
$ws.Range("E52").Value2 = 2
$ws.Range("F52").Value2 = "BTC"
$ws.Range("H52").Value2 = 0.0001
$ws.Range("I52").Value2 = "BTC"
$ws.Range("K52").Value2 = "Poloniex"
$ws.Range("L52").Value2 = 43093.501076388886

# Expand the table (ListObject) and its autofilter to cover the new row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:L61"))

# Expand the conditional-formatting ranges to cover the new last row.
$cf = $ws.Range("B2:B60").FormatConditions.Item(1)
$cf.ModifyAppliesToRange($ws.Range("B2:B61"))

$cf = $ws.Range("E2:E60").FormatConditions.Item(1)
$cf.ModifyAppliesToRange($ws.Range("E2:E61"))

$cf = $ws.Range("H2:H60").FormatConditions.Item(1)
$cf.ModifyAppliesToRange($ws.Range("H2:H61"))

# Restore the selection to A2 (top-left of the frozen/scrollable area).
$ws.Range("A2").Select()
